# Auto-generated edit script applying the cryptos.xlsx data refresh diff.
# Updates Price (D) and Volume(1h) (E) values for the updated rows, and
# for rows 12-14 also updates Coin (B) and Link (C) to reflect a re-ranking
# (WrappedliquidstakedEther2.0 / Polkadot / WrappedEther rows shifted order).
#
# Price values are stored as plain text (inline strings) in the source
# workbook, e.g. "214.48" rather than the number 214.48. Excel's COM layer
# auto-detects such single-dot numeric-looking strings and coerces them to
# real numbers on assignment, so for those cells we briefly force a text
# number format, assign the value, then restore the "Normal" style so the
# cell's style index matches the original (unstyled) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "25.999.68"
$ws.Range("E2").Value2 = "  -0.14%  "
$ws.Range("D3").Value2 = "1.631.28"
$ws.Range("E3").Value2 = "  -0.84%  "
$ws.Range("E4").Value2 = "  -0.18%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = "214.48"
$c.Style = "Normal"
$ws.Range("E5").Value2 = "  -0.75%  "
$ws.Range("E6").Value2 = "  -0.78%  "
$ws.Range("E7").Value2 = "  -0.18%  "
$ws.Range("E8").Value2 = "  -2.00%  "
$ws.Range("E9").Value2 = "  -3.07%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = "18.55"
$c.Style = "Normal"
$ws.Range("E10").Value2 = "  -5.32%  "
$ws.Range("E11").Value2 = "  -0.98%  "
$ws.Range("B12").Value2 = "WrappedEther"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value2 = "1.743.18"
$ws.Range("E12").Value2 = "  +5.67%  "
$ws.Range("B13").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value2 = "1.855.98"
$ws.Range("E13").Value2 = "  -0.94%  "
$ws.Range("B14").Value2 = "Polkadot"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = "4.19"
$c.Style = "Normal"
$ws.Range("E14").Value2 = "  -2.00%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = "0.529"
$c.Style = "Normal"
$ws.Range("E15").Value2 = "  -2.85%  "
$ws.Range("D16").Value2 = "26.010.46"
$ws.Range("E16").Value2 = "  -0.15%  "
$ws.Range("D17").Value2 = "0.0₃0742"
$ws.Range("E17").Value2 = "  -2.76%  "
$ws.Range("E18").Value2 = "  -3.07%  "
$ws.Range("E19").Value2 = "  -0.17%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = "193.50"
$c.Style = "Normal"
$ws.Range("E20").Value2 = "  -0.54%  "
$ws.Range("E21").Value2 = "  -2.40%  "
$ws.Range("E22").Value2 = "  -3.73%  "
$ws.Range("E23").Value2 = "  -1.99%  "
$ws.Range("E24").Value2 = "  +1.33%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = "144.07"
$c.Style = "Normal"
$ws.Range("E25").Value2 = "  +0.50%  "
$ws.Range("E26").Value2 = "  -0.06%  "
$ws.Range("E27").Value2 = "  -4.23%  "
$ws.Range("E28").Value2 = "  -1.89%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value2 = "15.29"
$c.Style = "Normal"
$ws.Range("E29").Value2 = "  -1.52%  "
$ws.Range("E30").Value2 = "  -1.03%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = "0.0485"
$c.Style = "Normal"
$ws.Range("E31").Value2 = "  -2.21%  "
$ws.Range("E32").Value2 = "  -3.92%  "
$ws.Range("E33").Value2 = "  -5.33%  "
$ws.Range("E35").Value2 = "  -2.58%  "
$ws.Range("D36").Value2 = "1.123.98"
$ws.Range("E36").Value2 = "  -0.68%  "
$ws.Range("E37").Value2 = "  -5.98%  "
$ws.Range("E38").Value2 = "  -1.48%  "
$ws.Range("E39").Value2 = "  -3.17%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = "0.0154"
$c.Style = "Normal"
$ws.Range("E40").Value2 = "  -2.13%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = "98.24"
$c.Style = "Normal"
$ws.Range("E41").Value2 = "  -0.79%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = "0.768"
$c.Style = "Normal"
$ws.Range("E42").Value2 = "  -3.67%  "
$ws.Range("D43").Value2 = "1.766.19"
$ws.Range("E43").Value2 = "  -0.92%  "
$ws.Range("E44").Value2 = "  -5.42%  "
$ws.Range("D45").Value2 = "0.0₆0114"
$ws.Range("E45").Value2 = "  -2.36%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = "54.58"
$c.Style = "Normal"
$ws.Range("E47").Value2 = "  -3.51%  "
$ws.Range("E48").Value2 = "  -0.52%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = "0.413"
$c.Style = "Normal"
$ws.Range("E49").Value2 = "  -0.55%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = "7.49"
$c.Style = "Normal"
$ws.Range("E50").Value2 = "  -3.81%  "
$ws.Range("E51").Value2 = "  -0.07%  "
